$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("numeric")

$ws.Range("B2").Value = "*"
$ws.Range("B3").Value = "*"

$ws.Range("B4").Select()
